# "excel file design finished"
# Adds the Right/Wrong/Not-Attempt/Max summary block (rows 9-12) and the
# two Student-Ans/Correct-Ans answer-key blocks (rows 15-40) below the
# existing "Mark Sheet" header, with Century-font / thin-border / centered
# styling and green/red/blue colour coding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colours as Excel COM BGR longs (matches Font.Color semantics).
$GREEN = 32768      # rgb 00008000
$BLUE  = 16711680   # rgb 000000FF
$RED   = 255        # rgb 00FF0000

# ---- Row 9: column headers for the summary block ----------------------
$ws.Cells.Item(9, 2).Value = "Right"
$ws.Cells.Item(9, 3).Value = "Wrong"
$ws.Cells.Item(9, 4).Value = "Not Attempt"
$ws.Cells.Item(9, 5).Value = "Max"

$row9 = $ws.Range("A9:E9")
$row9.Font.Name = "Century"
$row9.Font.Size = 12
$row9.Borders.LineStyle = 1
$row9.HorizontalAlignment = -4108

# ---- Rows 10-12: No./Marking/Total summary table -----------------------
$ws.Cells.Item(10, 1).Value = "No."
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = 12
$ws.Cells.Item(10, 4).Value = 11
$ws.Cells.Item(10, 5).Value = 28

$ws.Cells.Item(11, 1).Value = "Marking"
$ws.Cells.Item(11, 2).Value = 5
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 0

$ws.Cells.Item(12, 1).Value = "Total"
$ws.Cells.Item(12, 2).Value = 25
$ws.Cells.Item(12, 3).Value = -12
$ws.Cells.Item(12, 5).Value = "13/140"

$block = $ws.Range("A10:E12")
$block.Font.Name = "Century"
$block.Font.Size = 12
$block.Borders.LineStyle = 1
$block.HorizontalAlignment = -4108

# A/D/E columns of rows 10-12 stay plain (no colour override needed).
$greenCells = $ws.Range("B10:B12")
$greenCells.Font.Color = $GREEN

$redCells = $ws.Range("C10:C12")
$redCells.Font.Color = $RED

$ws.Cells.Item(12, 5).Font.Color = $BLUE

# ---- Row 15: "Student Ans" / "Correct Ans" headers (bold) --------------
$ws.Cells.Item(15, 1).Value = "Student Ans"
$ws.Cells.Item(15, 2).Value = "Correct Ans"
$ws.Cells.Item(15, 4).Value = "Student Ans"
$ws.Cells.Item(15, 5).Value = "Correct Ans"

$hdrLeft = $ws.Range("A15:B15")
$hdrLeft.Font.Name = "Century"
$hdrLeft.Font.Size = 12
$hdrLeft.Font.Bold = $true
$hdrLeft.Borders.LineStyle = 1
$hdrLeft.HorizontalAlignment = -4108

$hdrRight = $ws.Range("D15:E15")
$hdrRight.Font.Name = "Century"
$hdrRight.Font.Size = 12
$hdrRight.Font.Bold = $true
$hdrRight.Borders.LineStyle = 1
$hdrRight.HorizontalAlignment = -4108

# ---- Rows 16-40: answer key, columns A (student) / B (correct) --------
# "row|studentAns|correctAns"
$answersAB = @(
  "16|Option A|Option A",
  "17||Option D",
  "18|Option C|Option B",
  "19||Option C",
  "20||Option B",
  "21|Option D|Option C",
  "22|Option D|Option D",
  "23||Option D",
  "24||Option A",
  "25||Option A",
  "26|Option A|Option C",
  "27||Option A",
  "28|Option C|Option D",
  "29||Option D",
  "30|Option A|Option B",
  "31|Option A|Option D",
  "32|Option C|Option C",
  "33|Option C|Option D",
  "34||Option B",
  "35|Option C|Option D",
  "36|Option B|Option A",
  "37|Option A|Option A",
  "38|Option B|Option A",
  "39|Option B|Option D",
  "40||Option D"
)

# second answer key, columns D (student) / E (correct) -- rows 16-18 only
$answersDE = @(
  "16|Option A|Option A",
  "17||Option C",
  "18|Option A|Option D"
)

foreach ($line in $answersAB) {
    $parts = $line.Split("|")
    $r = [int]$parts[0]
    $stud = $parts[1]
    $corr = $parts[2]

    $studCell = $ws.Cells.Item($r, 1)
    if ($stud -ne "") {
        $studCell.Value = $stud
    }
    $studCell.Font.Name = "Century"
    $studCell.Font.Size = 12
    $studCell.Borders.LineStyle = 1
    $studCell.HorizontalAlignment = -4108
    if ($stud -eq $corr) {
        $studCell.Font.Color = $GREEN
    } else {
        $studCell.Font.Color = $RED
    }

    $corrCell = $ws.Cells.Item($r, 2)
    $corrCell.Value = $corr
    $corrCell.Font.Name = "Century"
    $corrCell.Font.Size = 12
    $corrCell.Borders.LineStyle = 1
    $corrCell.HorizontalAlignment = -4108
    $corrCell.Font.Color = $BLUE
}

foreach ($line in $answersDE) {
    $parts = $line.Split("|")
    $r = [int]$parts[0]
    $stud = $parts[1]
    $corr = $parts[2]

    $studCell = $ws.Cells.Item($r, 4)
    if ($stud -ne "") {
        $studCell.Value = $stud
    }
    $studCell.Font.Name = "Century"
    $studCell.Font.Size = 12
    $studCell.Borders.LineStyle = 1
    $studCell.HorizontalAlignment = -4108
    if ($stud -eq $corr) {
        $studCell.Font.Color = $GREEN
    } else {
        $studCell.Font.Color = $RED
    }

    $corrCell = $ws.Cells.Item($r, 5)
    $corrCell.Value = $corr
    $corrCell.Font.Name = "Century"
    $corrCell.Font.Size = 12
    $corrCell.Borders.LineStyle = 1
    $corrCell.HorizontalAlignment = -4108
    $corrCell.Font.Color = $BLUE
}

Write-Host "Applied mark-sheet design (rows 9-40)."
